# Generate Report for Handoff
#
# For the localization row corresponding to file
# "59df1511-dce7-46f6-9dc0-d61e9dbaecd1.md", a new handoff report run
# records fresh timestamps:
#   - zh-cn sheet: "Latest Handoff Datetime" (H7) -> 2016-09-04 05:02:01
#   - Overview sheet: "Latest HO Xliff Generate Date" (G7) -> 2016-09-04 05:02:10

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-04 05:02:01"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-04 05:02:10"
